$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert a new column before D (shifts existing D:K data to E:L)
$ws.Columns("D:D").Insert()

# Step 2: Copy number formats/fonts from column E (the old column D) into new column D
# (only for the row blocks that actually contain data, so we don't create stray
# cells/rows in the blank separator rows 5,6,36,37,78,79)
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# match column D width to column E for visual consistency
$ws.Columns("D:D").ColumnWidth = $ws.Columns("E:E").ColumnWidth

# Rows 7-35
$arr_7_35 = New-Object 'object[,]' 29,9
$arr_7_35[0,0] = 43465
$arr_7_35[0,1] = 43100
$arr_7_35[0,2] = 42735
$arr_7_35[0,3] = 42369
$arr_7_35[0,4] = 42004
$arr_7_35[0,5] = 41639
$arr_7_35[0,6] = 41274
$arr_7_35[0,7] = 40908
$arr_7_35[1,0] = 33171600
$arr_7_35[1,1] = 41846900
$arr_7_35[1,2] = 42829700
$arr_7_35[1,3] = 47859600
$arr_7_35[1,4] = 126891500
$arr_7_35[1,5] = 134288700
$arr_7_35[1,6] = 148207000
$arr_7_35[1,7] = 132586500
$arr_7_35[2,0] = 25135900
$arr_7_35[2,1] = 33035900
$arr_7_35[2,2] = 35665800
$arr_7_35[2,3] = 36653200
$arr_7_35[2,4] = 111786100
$arr_7_35[2,5] = 118231900
$arr_7_35[2,6] = 128852700
$arr_7_35[2,7] = 114239900
$arr_7_35[3,0] = 8035700
$arr_7_35[3,1] = 8811000
$arr_7_35[3,2] = 7163900
$arr_7_35[3,3] = 11206400
$arr_7_35[3,4] = 15105400
$arr_7_35[3,5] = 16056800
$arr_7_35[3,6] = 19354300
$arr_7_35[3,7] = 18346700
$arr_7_35[5,0] = "NA"
$arr_7_35[5,1] = "NA"
$arr_7_35[5,2] = "NA"
$arr_7_35[5,3] = "NA"
$arr_7_35[5,4] = 55000
$arr_7_35[5,5] = 79700
$arr_7_35[5,6] = 49400
$arr_7_35[5,7] = "NA"
$arr_7_35[6,0] = 0
$arr_7_35[6,1] = 0
$arr_7_35[6,2] = 0
$arr_7_35[6,3] = 0
$arr_7_35[6,4] = 0
$arr_7_35[6,5] = 0
$arr_7_35[6,6] = 0
$arr_7_35[6,7] = 0
$arr_7_35[7,0] = 33700
$arr_7_35[7,1] = 71800
$arr_7_35[7,2] = 106600
$arr_7_35[7,3] = 77400
$arr_7_35[7,4] = 102100
$arr_7_35[7,5] = 98700
$arr_7_35[7,6] = -320900
$arr_7_35[7,7] = 192500
$arr_7_35[8,0] = 1767100
$arr_7_35[8,1] = 1907400
$arr_7_35[8,2] = 4289400
$arr_7_35[8,3] = 6360600
$arr_7_35[8,4] = 9787100
$arr_7_35[8,5] = 5840000
$arr_7_35[8,6] = 5697500
$arr_7_35[8,7] = 8311700
$arr_7_35[10,0] = 28770100
$arr_7_35[10,1] = 36385000
$arr_7_35[10,2] = 43397500
$arr_7_35[10,3] = 47950500
$arr_7_35[10,4] = 127670100
$arr_7_35[10,5] = 128697900
$arr_7_35[10,6] = 143048100
$arr_7_35[10,7] = 133695800
$arr_7_35[11,0] = 4401600
$arr_7_35[11,1] = 5461800
$arr_7_35[11,2] = -567700
$arr_7_35[11,3] = -90900
$arr_7_35[11,4] = -778700
$arr_7_35[11,5] = 5590900
$arr_7_35[11,6] = 5158900
$arr_7_35[11,7] = -1109300
$arr_7_35[13,0] = -51600
$arr_7_35[13,1] = 908800
$arr_7_35[13,2] = -774200
$arr_7_35[13,3] = -617100
$arr_7_35[13,4] = -711300
$arr_7_35[13,5] = -825800
$arr_7_35[13,6] = 2199100
$arr_7_35[13,7] = 805200
$arr_7_35[14,0] = 6121300
$arr_7_35[14,1] = 8282600
$arr_7_35[14,2] = 2957800
$arr_7_35[14,3] = 5667900
$arr_7_35[14,4] = 8320700
$arr_7_35[14,5] = 10619100
$arr_7_35[14,6] = 13069200
$arr_7_35[14,7] = 8013400
$arr_7_35[15,0] = 665300
$arr_7_35[15,1] = 805600
$arr_7_35[15,2] = 593500
$arr_7_35[15,3] = 966000
$arr_7_35[15,4] = 1200500
$arr_7_35[15,5] = 1310500
$arr_7_35[15,6] = 3684600
$arr_7_35[15,7] = 3112900
$arr_7_35[16,0] = 3684600
$arr_7_35[16,1] = 5565100
$arr_7_35[16,2] = -1935400
$arr_7_35[16,3] = -1674000
$arr_7_35[16,4] = -2690500
$arr_7_35[16,5] = 3454600
$arr_7_35[16,6] = 3673400
$arr_7_35[16,7] = -3417000
$arr_7_35[17,0] = 51600
$arr_7_35[17,1] = 901000
$arr_7_35[17,2] = 493700
$arr_7_35[17,3] = 816800
$arr_7_35[17,4] = 639500
$arr_7_35[17,5] = 805600
$arr_7_35[17,6] = 783100
$arr_7_35[17,7] = -1216100
$arr_7_35[18,0] = 0
$arr_7_35[18,1] = 0
$arr_7_35[18,2] = 0
$arr_7_35[18,3] = 0
$arr_7_35[18,4] = 0
$arr_7_35[18,5] = 0
$arr_7_35[18,6] = 0
$arr_7_35[18,7] = 0
$arr_7_35[19,0] = 3633000
$arr_7_35[19,1] = 4664100
$arr_7_35[19,2] = -2429100
$arr_7_35[19,3] = -2490800
$arr_7_35[19,4] = -3330100
$arr_7_35[19,5] = 2649000
$arr_7_35[19,6] = 2890200
$arr_7_35[19,7] = -2200900
$arr_7_35[20,0] = 3295300
$arr_7_35[20,1] = 4378000
$arr_7_35[20,2] = 6049800
$arr_7_35[20,3] = -3188700
$arr_7_35[20,4] = -3363700
$arr_7_35[20,5] = 2236100
$arr_7_35[20,6] = 2414500
$arr_7_35[20,7] = -2621100
$arr_7_35[21,0] = 0
$arr_7_35[21,1] = 0
$arr_7_35[21,2] = 0
$arr_7_35[21,3] = 0
$arr_7_35[21,4] = 0
$arr_7_35[21,5] = 0
$arr_7_35[21,6] = 0
$arr_7_35[21,7] = 0
$arr_7_35[22,0] = 320900
$arr_7_35[22,1] = 25800
$arr_7_35[22,2] = -15530600
$arr_7_35[22,3] = -4664100
$arr_7_35[22,4] = -181800
$arr_7_35[22,5] = 110000
$arr_7_35[22,6] = 41500
$arr_7_35[22,7] = 16400
$arr_7_35[23,0] = 0
$arr_7_35[23,1] = 0
$arr_7_35[23,2] = 0
$arr_7_35[23,3] = 0
$arr_7_35[23,4] = 0
$arr_7_35[23,5] = 0
$arr_7_35[23,6] = 0
$arr_7_35[23,7] = 0
$arr_7_35[24,0] = 0
$arr_7_35[24,1] = 0
$arr_7_35[24,2] = 0
$arr_7_35[24,3] = 0
$arr_7_35[24,4] = 0
$arr_7_35[24,5] = 0
$arr_7_35[24,6] = 0
$arr_7_35[24,7] = 0
$arr_7_35[25,0] = 51600
$arr_7_35[25,1] = -908800
$arr_7_35[25,2] = 774200
$arr_7_35[25,3] = 617100
$arr_7_35[25,4] = 711300
$arr_7_35[25,5] = 825800
$arr_7_35[25,6] = -2199100
$arr_7_35[25,7] = -805200
$arr_7_35[26,0] = 3616200
$arr_7_35[26,1] = 4403800
$arr_7_35[26,2] = -9480800
$arr_7_35[26,3] = -7852800
$arr_7_35[26,4] = -3545500
$arr_7_35[26,5] = 2346100
$arr_7_35[26,6] = 2456000
$arr_7_35[26,7] = -2604700
$arr_7_35[27,0] = 0
$arr_7_35[27,1] = 0
$arr_7_35[27,2] = 0
$arr_7_35[27,3] = 0
$arr_7_35[27,4] = 0
$arr_7_35[27,5] = 0
$arr_7_35[27,6] = 0
$arr_7_35[27,7] = 0
$arr_7_35[28,0] = 3616200
$arr_7_35[28,1] = 4403800
$arr_7_35[28,2] = -9480800
$arr_7_35[28,3] = -7852800
$arr_7_35[28,4] = -3545500
$arr_7_35[28,5] = 2346100
$arr_7_35[28,6] = 2456000
$arr_7_35[28,7] = -2604700
$ws.Range("D7:L35").Value = $arr_7_35

# Rows 38-77
$arr_38_77 = New-Object 'object[,]' 40,9
$arr_38_77[0,0] = 43465
$arr_38_77[0,1] = 43100
$arr_38_77[0,2] = 42735
$arr_38_77[0,3] = 42369
$arr_38_77[0,4] = 42004
$arr_38_77[0,5] = 41639
$arr_38_77[0,6] = 41274
$arr_38_77[0,7] = 40908
$arr_38_77[3,0] = 4402700
$arr_38_77[3,1] = 3038300
$arr_38_77[3,2] = 6254000
$arr_38_77[3,3] = 5822000
$arr_38_77[3,4] = 3580300
$arr_38_77[3,5] = 5079200
$arr_38_77[3,6] = 3159500
$arr_38_77[3,7] = 4521500
$arr_38_77[4,0] = 868400
$arr_38_77[4,1] = 751700
$arr_38_77[4,2] = 2408900
$arr_38_77[4,3] = 2331500
$arr_38_77[4,4] = 2033000
$arr_38_77[4,5] = 2971000
$arr_38_77[4,6] = 7923500
$arr_38_77[4,7] = 3614200
$arr_38_77[5,0] = 4950200
$arr_38_77[5,1] = 5193700
$arr_38_77[5,2] = 9012900
$arr_38_77[5,3] = 31588500
$arr_38_77[5,4] = 30778400
$arr_38_77[5,5] = 51446600
$arr_38_77[5,6] = 85039000
$arr_38_77[5,7] = 40375500
$arr_38_77[6,0] = 767400
$arr_38_77[6,1] = 890900
$arr_38_77[6,2] = 880800
$arr_38_77[6,3] = 2856600
$arr_38_77[6,4] = 3765400
$arr_38_77[6,5] = 9305800
$arr_38_77[6,6] = 10624100
$arr_38_77[6,7] = 5667200
$arr_38_77[7,0] = 15311800
$arr_38_77[7,1] = 7837100
$arr_38_77[7,2] = 969400
$arr_38_77[7,3] = 2371900
$arr_38_77[7,4] = 7667700
$arr_38_77[7,5] = 1873700
$arr_38_77[7,6] = 10626400
$arr_38_77[7,7] = 5276300
$arr_38_77[8,0] = 26300600
$arr_38_77[8,1] = 17711700
$arr_38_77[8,2] = 19526000
$arr_38_77[8,3] = 44970500
$arr_38_77[8,4] = 47824800
$arr_38_77[8,5] = 41233100
$arr_38_77[8,6] = 50399800
$arr_38_77[8,7] = 59454700
$arr_38_77[9,0] = 6673600
$arr_38_77[9,1] = 8459800
$arr_38_77[9,2] = 15507000
$arr_38_77[9,3] = 22005600
$arr_38_77[9,4] = 21234800
$arr_38_77[9,5] = 21158500
$arr_38_77[9,6] = 31260900
$arr_38_77[9,7] = 23176900
$arr_38_77[10,0] = 20259800
$arr_38_77[10,1] = 27787200
$arr_38_77[10,2] = 28321300
$arr_38_77[10,3] = 43754200
$arr_38_77[10,4] = 46307900
$arr_38_77[10,5] = 112385000
$arr_38_77[10,6] = 5944300
$arr_38_77[10,7] = 65579600
$arr_38_77[11,0] = 4730300
$arr_38_77[11,1] = 6260700
$arr_38_77[11,2] = 6498600
$arr_38_77[11,3] = 12236400
$arr_38_77[11,4] = 18730500
$arr_38_77[11,5] = 43340200
$arr_38_77[11,6] = 58310900
$arr_38_77[11,7] = 25184100
$arr_38_77[12,0] = 0
$arr_38_77[12,1] = 0
$arr_38_77[12,2] = 0
$arr_38_77[12,3] = 0
$arr_38_77[12,4] = 0
$arr_38_77[12,5] = 0
$arr_38_77[12,6] = 0
$arr_38_77[12,7] = 0
$arr_38_77[13,0] = 0
$arr_38_77[13,1] = 0
$arr_38_77[13,2] = 0
$arr_38_77[13,3] = 0
$arr_38_77[13,4] = 0
$arr_38_77[13,5] = 0
$arr_38_77[13,6] = 0
$arr_38_77[13,7] = 0
$arr_38_77[14,0] = 2986700
$arr_38_77[14,1] = 2555900
$arr_38_77[14,2] = 1616800
$arr_38_77[14,3] = 4595700
$arr_38_77[14,4] = 6924900
$arr_38_77[14,5] = 8218600
$arr_38_77[14,6] = 6150700
$arr_38_77[14,7] = 6047500
$arr_38_77[15,0] = 0
$arr_38_77[15,1] = 0
$arr_38_77[15,2] = 0
$arr_38_77[15,3] = 0
$arr_38_77[15,4] = 0
$arr_38_77[15,5] = 0
$arr_38_77[15,6] = 0
$arr_38_77[15,7] = 0
$arr_38_77[16,0] = 60951000
$arr_38_77[16,1] = 62775300
$arr_38_77[16,2] = 71469600
$arr_38_77[16,3] = 127562000
$arr_38_77[16,4] = 141023000
$arr_38_77[16,5] = 148473000
$arr_38_77[16,6] = 159642000
$arr_38_77[16,7] = 179443000
$arr_38_77[19,0] = 1862500
$arr_38_77[19,1] = 2019600
$arr_38_77[19,2] = 2288900
$arr_38_77[19,3] = 2664700
$arr_38_77[19,4] = 2451500
$arr_38_77[19,5] = 2788100
$arr_38_77[19,6] = 6124900
$arr_38_77[19,7] = "NA"
$arr_38_77[20,0] = 1753700
$arr_38_77[20,1] = 3477000
$arr_38_77[20,2] = 4254600
$arr_38_77[20,3] = 3128100
$arr_38_77[20,4] = 4356700
$arr_38_77[20,5] = 5243100
$arr_38_77[20,6] = 8557400
$arr_38_77[20,7] = 6907900
$arr_38_77[21,0] = 13506500
$arr_38_77[21,1] = 10260600
$arr_38_77[21,2] = 19402600
$arr_38_77[21,3] = 31731000
$arr_38_77[21,4] = 33181700
$arr_38_77[21,5] = 52522600
$arr_38_77[21,6] = 86580600
$arr_38_77[21,7] = 47240000
$arr_38_77[22,0] = 17122700
$arr_38_77[22,1] = 15757200
$arr_38_77[22,2] = 25946000
$arr_38_77[22,3] = 37523800
$arr_38_77[22,4] = 39990000
$arr_38_77[22,5] = 36479300
$arr_38_77[22,6] = 40464600
$arr_38_77[22,7] = 54147900
$arr_38_77[23,0] = 9338300
$arr_38_77[23,1] = 11132400
$arr_38_77[23,2] = 11708000
$arr_38_77[23,3] = 16778200
$arr_38_77[23,4] = 17709500
$arr_38_77[23,5] = 20253000
$arr_38_77[23,6] = 24421200
$arr_38_77[23,7] = 28205500
$arr_38_77[24,0] = 24932900
$arr_38_77[24,1] = 28359400
$arr_38_77[24,2] = 32371700
$arr_38_77[24,3] = 51856100
$arr_38_77[24,4] = 53351700
$arr_38_77[24,5] = 85963500
$arr_38_77[24,6] = 118969000
$arr_38_77[24,7] = 50591200
$arr_38_77[25,0] = 0
$arr_38_77[25,1] = 0
$arr_38_77[25,2] = 0
$arr_38_77[25,3] = 0
$arr_38_77[25,4] = 0
$arr_38_77[25,5] = 0
$arr_38_77[25,6] = 0
$arr_38_77[25,7] = 0
$arr_38_77[26,0] = 0
$arr_38_77[26,1] = 0
$arr_38_77[26,2] = 0
$arr_38_77[26,3] = 0
$arr_38_77[26,4] = 0
$arr_38_77[26,5] = 0
$arr_38_77[26,6] = 0
$arr_38_77[26,7] = 0
$arr_38_77[27,0] = 0
$arr_38_77[27,1] = 0
$arr_38_77[27,2] = 0
$arr_38_77[27,3] = 0
$arr_38_77[27,4] = 0
$arr_38_77[27,5] = 0
$arr_38_77[27,6] = 0
$arr_38_77[27,7] = 0
$arr_38_77[28,0] = 54490600
$arr_38_77[28,1] = 58279500
$arr_38_77[28,2] = 72653300
$arr_38_77[28,3] = 109129000
$arr_38_77[28,4] = 113439000
$arr_38_77[28,5] = 110636000
$arr_38_77[28,6] = 120079000
$arr_38_77[28,7] = 137494000
$arr_38_77[30,0] = 0
$arr_38_77[30,1] = 0
$arr_38_77[30,2] = 0
$arr_38_77[30,3] = 0
$arr_38_77[30,4] = 0
$arr_38_77[30,5] = 0
$arr_38_77[30,6] = 0
$arr_38_77[30,7] = 0
$arr_38_77[31,0] = 0
$arr_38_77[31,1] = 0
$arr_38_77[31,2] = 0
$arr_38_77[31,3] = 0
$arr_38_77[31,4] = 0
$arr_38_77[31,5] = 0
$arr_38_77[31,6] = 0
$arr_38_77[31,7] = 0
$arr_38_77[32,0] = 0
$arr_38_77[32,1] = 0
$arr_38_77[32,2] = 0
$arr_38_77[32,3] = 0
$arr_38_77[32,4] = 0
$arr_38_77[32,5] = 0
$arr_38_77[32,6] = 0
$arr_38_77[32,7] = 0
$arr_38_77[33,0] = 0
$arr_38_77[33,1] = 0
$arr_38_77[33,2] = 0
$arr_38_77[33,3] = 0
$arr_38_77[33,4] = 0
$arr_38_77[33,5] = 0
$arr_38_77[33,6] = 0
$arr_38_77[33,7] = 0
$arr_38_77[34,0] = -2761200
$arr_38_77[34,1] = -5107300
$arr_38_77[34,2] = -9531300
$arr_38_77[34,3] = 10568000
$arr_38_77[34,4] = 18896600
$arr_38_77[34,5] = 26149100
$arr_38_77[34,6] = 26908700
$arr_38_77[34,7] = 45118900
$arr_38_77[35,0] = 0
$arr_38_77[35,1] = 0
$arr_38_77[35,2] = 0
$arr_38_77[35,3] = 0
$arr_38_77[35,4] = 0
$arr_38_77[35,5] = 0
$arr_38_77[35,6] = 0
$arr_38_77[35,7] = 0
$arr_38_77[36,0] = 0
$arr_38_77[36,1] = 0
$arr_38_77[36,2] = 0
$arr_38_77[36,3] = 0
$arr_38_77[36,4] = 0
$arr_38_77[36,5] = 0
$arr_38_77[36,6] = 0
$arr_38_77[36,7] = 0
$arr_38_77[37,0] = 0
$arr_38_77[37,1] = 0
$arr_38_77[37,2] = 0
$arr_38_77[37,3] = 0
$arr_38_77[37,4] = 0
$arr_38_77[37,5] = 0
$arr_38_77[37,6] = 0
$arr_38_77[37,7] = 0
$arr_38_77[38,0] = 6460400
$arr_38_77[38,1] = 4495800
$arr_38_77[38,2] = -1183700
$arr_38_77[38,3] = 18433200
$arr_38_77[38,4] = 27584100
$arr_38_77[38,5] = 37836900
$arr_38_77[38,6] = 39563600
$arr_38_77[38,7] = 41948400
$arr_38_77[39,0] = 0
$arr_38_77[39,1] = 0
$arr_38_77[39,2] = 0
$arr_38_77[39,3] = 0
$arr_38_77[39,4] = 0
$arr_38_77[39,5] = 0
$arr_38_77[39,6] = 0
$arr_38_77[39,7] = 0
$ws.Range("D38:L77").Value = $arr_38_77

# Rows 80-102
$arr_80_102 = New-Object 'object[,]' 23,9
$arr_80_102[0,0] = 43465
$arr_80_102[0,1] = 43100
$arr_80_102[0,2] = 42735
$arr_80_102[0,3] = 42369
$arr_80_102[0,4] = 42004
$arr_80_102[0,5] = 41639
$arr_80_102[0,6] = 41274
$arr_80_102[0,7] = 40908
$arr_80_102[1,0] = 3616200
$arr_80_102[1,1] = 4403800
$arr_80_102[1,2] = -9480800
$arr_80_102[1,3] = -7852800
$arr_80_102[1,4] = -3545500
$arr_80_102[1,5] = 2346100
$arr_80_102[1,6] = 2456000
$arr_80_102[1,7] = -2604700
$arr_80_102[3,0] = 1767100
$arr_80_102[3,1] = 1907400
$arr_80_102[3,2] = 4289400
$arr_80_102[3,3] = 6360600
$arr_80_102[3,4] = 9787100
$arr_80_102[3,5] = 5840000
$arr_80_102[3,6] = 5697500
$arr_80_102[3,7] = 8311700
$arr_80_102[4,0] = 0
$arr_80_102[4,1] = 0
$arr_80_102[4,2] = 0
$arr_80_102[4,3] = 0
$arr_80_102[4,4] = 0
$arr_80_102[4,5] = 0
$arr_80_102[4,6] = 0
$arr_80_102[4,7] = 0
$arr_80_102[5,0] = 0
$arr_80_102[5,1] = 0
$arr_80_102[5,2] = 0
$arr_80_102[5,3] = 0
$arr_80_102[5,4] = 0
$arr_80_102[5,5] = 0
$arr_80_102[5,6] = 0
$arr_80_102[5,7] = 0
$arr_80_102[6,0] = 0
$arr_80_102[6,1] = 0
$arr_80_102[6,2] = 0
$arr_80_102[6,3] = 0
$arr_80_102[6,4] = 0
$arr_80_102[6,5] = 0
$arr_80_102[6,6] = 0
$arr_80_102[6,7] = 0
$arr_80_102[7,0] = 0
$arr_80_102[7,1] = 0
$arr_80_102[7,2] = 0
$arr_80_102[7,3] = 0
$arr_80_102[7,4] = 0
$arr_80_102[7,5] = 0
$arr_80_102[7,6] = 0
$arr_80_102[7,7] = 0
$arr_80_102[8,0] = 0
$arr_80_102[8,1] = 0
$arr_80_102[8,2] = 0
$arr_80_102[8,3] = 0
$arr_80_102[8,4] = 0
$arr_80_102[8,5] = 0
$arr_80_102[8,6] = 0
$arr_80_102[8,7] = 0
$arr_80_102[9,0] = 3127000
$arr_80_102[9,1] = -3413100
$arr_80_102[9,2] = 5750200
$arr_80_102[9,3] = 6932800
$arr_80_102[9,4] = 7268300
$arr_80_102[9,5] = 7235700
$arr_80_102[9,6] = 9882500
$arr_80_102[9,7] = 7758900
$arr_80_102[11,0] = -2558100
$arr_80_102[11,1] = -2301200
$arr_80_102[11,2] = -3405200
$arr_80_102[11,3] = -3345800
$arr_80_102[11,4] = -4484600
$arr_80_102[11,5] = -5026500
$arr_80_102[11,6] = -7157200
$arr_80_102[11,7] = -7296400
$arr_80_102[12,0] = 0
$arr_80_102[12,1] = 0
$arr_80_102[12,2] = 0
$arr_80_102[12,3] = 0
$arr_80_102[12,4] = 0
$arr_80_102[12,5] = 0
$arr_80_102[12,6] = 0
$arr_80_102[12,7] = 0
$arr_80_102[13,0] = 0
$arr_80_102[13,1] = 0
$arr_80_102[13,2] = 0
$arr_80_102[13,3] = 0
$arr_80_102[13,4] = 0
$arr_80_102[13,5] = 0
$arr_80_102[13,6] = 0
$arr_80_102[13,7] = 0
$arr_80_102[14,0] = 1134300
$arr_80_102[14,1] = -438700
$arr_80_102[14,2] = -4898600
$arr_80_102[14,3] = -322000
$arr_80_102[14,4] = -3629600
$arr_80_102[14,5] = -765200
$arr_80_102[14,6] = -3386200
$arr_80_102[14,7] = -3581300
$arr_80_102[16,0] = -729300
$arr_80_102[16,1] = -387100
$arr_80_102[16,2] = -1095100
$arr_80_102[16,3] = -792100
$arr_80_102[16,4] = -942500
$arr_80_102[16,5] = -2352800
$arr_80_102[16,6] = -2137400
$arr_80_102[16,7] = -3354700
$arr_80_102[17,0] = 0
$arr_80_102[17,1] = 0
$arr_80_102[17,2] = 0
$arr_80_102[17,3] = 0
$arr_80_102[17,4] = 0
$arr_80_102[17,5] = 0
$arr_80_102[17,6] = 0
$arr_80_102[17,7] = 0
$arr_80_102[18,0] = 0
$arr_80_102[18,1] = 0
$arr_80_102[18,2] = 0
$arr_80_102[18,3] = 0
$arr_80_102[18,4] = 0
$arr_80_102[18,5] = 0
$arr_80_102[18,6] = 0
$arr_80_102[18,7] = 0
$arr_80_102[19,0] = 0
$arr_80_102[19,1] = 0
$arr_80_102[19,2] = 0
$arr_80_102[19,3] = 0
$arr_80_102[19,4] = 0
$arr_80_102[19,5] = 0
$arr_80_102[19,6] = 0
$arr_80_102[19,7] = 0
$arr_80_102[20,0] = -2958700
$arr_80_102[20,1] = 605900
$arr_80_102[20,2] = -323100
$arr_80_102[20,3] = -4328600
$arr_80_102[20,4] = -5173500
$arr_80_102[20,5] = -4479000
$arr_80_102[20,6] = -7683400
$arr_80_102[20,7] = -6849200
$arr_80_102[21,0] = 0
$arr_80_102[21,1] = -9000
$arr_80_102[21,2] = -97600
$arr_80_102[21,3] = -67300
$arr_80_102[21,4] = 50500
$arr_80_102[21,5] = -66200
$arr_80_102[21,6] = 29200
$arr_80_102[21,7] = -14100
$arr_80_102[22,0] = 1302600
$arr_80_102[22,1] = -3254900
$arr_80_102[22,2] = 430800
$arr_80_102[22,3] = 2214800
$arr_80_102[22,4] = -1484400
$arr_80_102[22,5] = 1925300
$arr_80_102[22,6] = -1157900
$arr_80_102[22,7] = -2685700
$ws.Range("D80:L102").Value = $arr_80_102
